$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2775001
$ws.Range("I62").Value = 3533334.8
$ws.Range("K62").Value = 3533334.8
$ws.Range("M62").Value = -3532710.8
$ws.Range("H65").Value = 2775001
$ws.Range("I65").Value = 3533334.8
$ws.Range("K65").Value = 17666674
$ws.Range("M65").Value = -17663554
$ws.Range("H98").Value = 1609.5625
$ws.Range("I98").Value = 1520.6
$ws.Range("K98").Value = 1520.6
$ws.Range("M98").Value = -22.59999999999991
$ws.Range("H106").Value = 262500
$ws.Range("J106").Value = 505000
$ws.Range("L106").Value = 505000
$ws.Range("N106").Value = -506262
$ws.Range("H122").Value = 1609.5625
$ws.Range("I122").Value = 1520.6
$ws.Range("K122").Value = 4561.799999999999
$ws.Range("M122").Value = -2111.799999999999
$ws.Range("H129").Value = 1553.3334
$ws.Range("I129").Value = 1211.5714
$ws.Range("J129").Value = 2749.5
$ws.Range("K129").Value = 3634.7142
$ws.Range("L129").Value = 8248.5
$ws.Range("M129").Value = 1365.2858
$ws.Range("N129").Value = -18248.5
$ws.Range("H132").Value = 4250.844
$ws.Range("I132").Value = 3839.3513
$ws.Range("K132").Value = 11518.0539
$ws.Range("M132").Value = -8988.053899999999
$ws.Range("H133").Value = 72779.336
$ws.Range("J133").Value = 72779.336
$ws.Range("L133").Value = 72779.336
$ws.Range("N133").Value = -82899.336
$ws.Range("H134").Value = 75708
$ws.Range("J134").Value = 75708
$ws.Range("L134").Value = 75708
$ws.Range("N134").Value = -85848
$ws.Range("H137").Value = 2199.8823
$ws.Range("I137").Value = 1947.84
$ws.Range("J137").Value = 2900
$ws.Range("K137").Value = 5843.52
$ws.Range("L137").Value = 8700
$ws.Range("M137").Value = -3293.52
$ws.Range("N137").Value = -13800
$ws.Range("H138").Value = 2406.2563
$ws.Range("I138").Value = 1420
$ws.Range("K138").Value = 4260
$ws.Range("M138").Value = 880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 140147.36
$ws.Range("I32").Value = 149534.73
$ws.Range("K32").Value = 149534.73
$ws.Range("M32").Value = -149247.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3866.4546
$ws.Range("I94").Value = 3614.5557
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 3614.5557
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -3163.5557
$ws.Range("N94").Value = -5902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2382709
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 2647232.2
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2647232.2
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -2647456.2
$ws.Range("H6").Value = 814188.1
$ws.Range("I6").Value = 907776.75
$ws.Range("J6").Value = 34283
$ws.Range("K6").Value = 907776.75
$ws.Range("L6").Value = 34283
$ws.Range("M6").Value = -907663.75
$ws.Range("N6").Value = -34509
$ws.Range("H22").Value = 301.15625
$ws.Range("I22").Value = 337.1579
$ws.Range("K22").Value = 337.1579
$ws.Range("M22").Value = 12.84210000000002
$ws.Range("H31").Value = 4688.6567
$ws.Range("I31").Value = 9993.929
$ws.Range("J31").Value = 3287.2642
$ws.Range("K31").Value = 9993.929
$ws.Range("L31").Value = 3287.2642
$ws.Range("M31").Value = -9698.929
$ws.Range("N31").Value = -3877.2642
$ws.Range("H34").Value = 4688.6567
$ws.Range("I34").Value = 9993.929
$ws.Range("J34").Value = 3287.2642
$ws.Range("K34").Value = 9993.929
$ws.Range("L34").Value = 3287.2642
$ws.Range("M34").Value = -9791.929
$ws.Range("N34").Value = -3691.2642
$ws.Range("H56").Value = 16000
$ws.Range("I56").Value = 16000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 16000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -15155
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 835769.56
$ws.Range("I58").Value = 2164.8333
$ws.Range("K58").Value = 2164.8333
$ws.Range("M58").Value = -1961.8333
$ws.Range("H62").Value = 8333
$ws.Range("J62").Value = 7499
$ws.Range("L62").Value = 7499
$ws.Range("N62").Value = -8747
$ws.Range("H65").Value = 8333
$ws.Range("J65").Value = 7499
$ws.Range("L65").Value = 37495
$ws.Range("N65").Value = -43735
$ws.Range("H94").Value = 1815.5
$ws.Range("I94").Value = 1340.75
$ws.Range("J94").Value = 2132
$ws.Range("K94").Value = 1340.75
$ws.Range("L94").Value = 2132
$ws.Range("M94").Value = -889.75
$ws.Range("N94").Value = -3034
$ws.Range("H136").Value = 835769.56
$ws.Range("I136").Value = 2164.8333
$ws.Range("K136").Value = 6494.499899999999
$ws.Range("M136").Value = -3944.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132.45454
$ws.Range("I2").Value = 74.666664
$ws.Range("J2").Value = 154.125
$ws.Range("K2").Value = 447.999984
$ws.Range("L2").Value = 924.75
$ws.Range("M2").Value = -334.999984
$ws.Range("N2").Value = -1150.75
$ws.Range("H101").Value = 2000
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -10868
$ws.Range("H122").Value = 664.1667
$ws.Range("I122").Value = 750.53845
$ws.Range("K122").Value = 6754.84605
$ws.Range("M122").Value = -4304.84605
$ws.Range("H131").Value = 5791.857
$ws.Range("I131").Value = 1843
$ws.Range("K131").Value = 5529
$ws.Range("M131").Value = -489

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 83750
$ws.Range("I5").Value = 83750
$ws.Range("K5").Value = 83750
$ws.Range("M5").Value = -83638
$ws.Range("H44").Value = 26676
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -31192
$ws.Range("H52").Value = 30581.334
$ws.Range("I52").Value = 28374.5
$ws.Range("K52").Value = 28374.5
$ws.Range("M52").Value = -28115.5
$ws.Range("H122").Value = 349503.12
$ws.Range("I122").Value = 429273.06
$ws.Range("K122").Value = 1287819.18
$ws.Range("M122").Value = -1285369.18
$ws.Range("H126").Value = 8511.667
$ws.Range("I126").Value = 9075.625
$ws.Range("K126").Value = 27226.875
$ws.Range("M126").Value = -24756.875
$ws.Range("H132").Value = 276223.72
$ws.Range("I132").Value = 327983.38
$ws.Range("K132").Value = 983950.14
$ws.Range("M132").Value = -981420.14
$ws.Range("H138").Value = 61900
$ws.Range("J138").Value = 61900
$ws.Range("L138").Value = 61900
$ws.Range("N138").Value = -72180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 500900
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 715142.9
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 715142.9
$ws.Range("N2").Value = -715366.9
$ws.Range("M2").Value = -888
$ws.Range("H41").Value = 10031
$ws.Range("I41").Value = 10031
$ws.Range("K41").Value = 10031
$ws.Range("M41").Value = -9593
$ws.Range("H68").Value = 7179.8335
$ws.Range("J68").Value = 8050
$ws.Range("L68").Value = 8050
$ws.Range("N68").Value = -9548
$ws.Range("H71").Value = 7179.8335
$ws.Range("J71").Value = 8050
$ws.Range("L71").Value = 40250
$ws.Range("N71").Value = -47738
$ws.Range("H93").Value = 6787.25
$ws.Range("I93").Value = 4500
$ws.Range("K93").Value = 4500
$ws.Range("M93").Value = -3252
$ws.Range("H105").Value = 247749.75
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 247749.75
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 247749.75
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -254737.75
$ws.Range("H120").Value = 89999
$ws.Range("J120").Value = 89999
$ws.Range("L120").Value = 89999
$ws.Range("N120").Value = -99675
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H136").Value = 2240.3242
$ws.Range("I136").Value = 2006.0605
$ws.Range("K136").Value = 6018.181500000001
$ws.Range("M136").Value = -3468.181500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 901
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 1500
$ws.Range("N2").Value = -1724
$ws.Range("H18").Value = 23949.5
$ws.Range("J18").Value = 23949.5
$ws.Range("L18").Value = 23949.5
$ws.Range("N18").Value = -24295.5
$ws.Range("H104").Value = 12194.833
$ws.Range("J104").Value = 12194.833
$ws.Range("L104").Value = 12194.833
$ws.Range("N104").Value = -19182.833
$ws.Range("H126").Value = 2031.6666
$ws.Range("I126").Value = 1760.125
$ws.Range("K126").Value = 5280.375
$ws.Range("M126").Value = -2810.375
